# daily auto push: 2026-01-30 19:00 UTC
#
# A new reading came in for 2026/01/30 (金) at hour 23. Chronologically it
# belongs right after the existing 2026/01/30 block (rows 738-741) and
# before the 2026/12/29 block, which currently starts at row 742. So a
# new row is inserted at 742, pushing every row below it down by one; the
# sheet's used range grows from A1:D783 to A1:D784 as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$insertRow = 742

# Push row 742 (and everything after it) down by one row.
$ws.Rows.Item($insertRow).Insert()

# Column A stores the date as plain text (e.g. "2026/01/30"), not a real
# Excel date serial, and column B stores the weekday kanji as plain text
# too. Assigning those literal strings straight to .Value would make
# Excel auto-detect/convert "2026/01/30" into a date (stamping a date
# number format onto the cell), which the source data doesn't have.
# Instead, copy the already-correct text+format from the existing
# 2026/01/30 rows right above the insertion point - Range.Copy carries
# over the resolved text value verbatim, with no re-parsing, and without
# touching the shared style table.
$ws.Cells.Item(738, 1).Copy($ws.Cells.Item($insertRow, 1))
$ws.Cells.Item(738, 2).Copy($ws.Cells.Item($insertRow, 2))

$ws.Cells.Item($insertRow, 3).Value = 23
$ws.Cells.Item($insertRow, 4).Value = 201
